$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (trial numbers) for columns B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update CON row (row 2) values for columns B2:E2
$ws.Range("B2").Value = 108.29594798993637
$ws.Range("C2").Value = 105.98916852820224
$ws.Range("D2").Value = 107.2116753546531
$ws.Range("E2").Value = 107.76606483851549

# Update STR row (row 3) values for columns B3:E3
$ws.Range("B3").Value = 107.03031794451725
$ws.Range("C3").Value = 104.88524901633632
$ws.Range("D3").Value = 107.2281202662675
$ws.Range("E3").Value = 108.64319819792583

# Update selection to match the new selected range used when the workbook was last saved
$ws.Range("B1:E3").Select()
